$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: new "Électroaimant+shipping" purchase line -------------------
# C10: article name (new shared string)
$ws.Range("C10").Value = "Électroaimant+shipping"
# D10: quantity
$ws.Range("D10").Value = 1
# E10: unit price, computed formula (12.79 + 6.48 shipping), default (no) style
$ws.Range("E10").Formula = "=12.79+6.48"
$ws.Range("E10").Style = "Normal"
# G10: purchased-by ("MS", reuses existing shared string)
$ws.Range("G10").Value = "MS"
# J10: amount already reimbursed (literal value, not a formula)
$ws.Range("J10").Value = 12.79

# --- Row 11: give J11 its own currency number format + distinct border ----
# (mirrors an existing "thin left/right + hair bottom" bordered, currency
#  formatted cell elsewhere on the sheet, minus its top border)
$ws.Range("C11").Copy() | Out-Null
$ws.Range("J11").PasteSpecial(-4122) | Out-Null
$ws.Range("J11").Borders.Item(8).LineStyle = -4142
$ws.Range("J11").NumberFormat = "#,##0.00\ ""$"""

# --- L7: roll the new H10 (shipping share) into the MS subtotal -----------
$ws.Range("L7").Formula = "=H4+H5+H10"

# --- cosmetic: restore the last-saved selection ----------------------------
$ws.Range("L18").Select() | Out-Null
